$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7967572212219238
$ws.Range("B1").Value = 3.103760480880737
$ws.Range("C1").Value = 2.955079555511475
$ws.Range("D1").Value = 2.510106086730957
$ws.Range("E1").Value = 2.154290676116943
